$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '77.035.39'
$ws.Range("E2").Value = '  +0.65%  '
# Row 3
$ws.Range("D3").Value = '2.953.36'
$ws.Range("E3").Value = '  +2.42%  '
# Row 4
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").Value = "'200.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.53%  '
# Row 6
$ws.Range("D6").Value = "'598.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '
# Row 7
$ws.Range("E7").Value = '  +0.00%  '
# Row 8
$ws.Range("D8").Value = "'0.551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.67%  '
# Row 9
$ws.Range("E9").Value = '  +1.76%  '
# Row 10
$ws.Range("D10").Value = '2.950.79'
$ws.Range("E10").Value = '  +2.41%  '
# Row 11
$ws.Range("D11").Value = "'0.443"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +12.23%  '
# Row 12
$ws.Range("E12").Value = '  +0.45%  '
# Row 13
$ws.Range("D13").Value = '3.496.18'
$ws.Range("E13").Value = '  +2.45%  '
# Row 14
$ws.Range("E14").Value = '  -0.67%  '
# Row 15
$ws.Range("D15").Value = '76.931.28'
$ws.Range("E15").Value = '  +0.68%  '
# Row 16
$ws.Range("D16").Value = "'28.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.04%  '
# Row 17
$ws.Range("D17").Value = "'0.0000189"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.40%  '
# Row 18
$ws.Range("D18").Value = '2.963.12'
$ws.Range("E18").Value = '  +3.08%  '
# Row 19
$ws.Range("D19").Value = "'13.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.75%  '
# Row 20
$ws.Range("E20").Value = '  -3.94%  '
# Row 21
$ws.Range("D21").Value = "'375.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.95%  '
# Row 22
$ws.Range("E22").Value = '  +4.85%  '
# Row 23
$ws.Range("D23").Value = "'2.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.34%  '
# Row 24
$ws.Range("D24").Value = "'72.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.07%  '
# Row 25
$ws.Range("E25").Value = '  +2.85%  '
# Row 26
$ws.Range("E26").Value = '  -0.07%  '
# Row 27
$ws.Range("D27").Value = "'4.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.46%  '
# Row 28
$ws.Range("D28").Value = "'9.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.66%  '
# Row 29
$ws.Range("E29").Value = '  +1.70%  '
# Row 30
$ws.Range("E30").Value = '  +0.25%  '
# Row 31
$ws.Range("D31").Value = "'8.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.19%  '
# Row 32
$ws.Range("D32").Value = "'1.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.21%  '
# Row 33
$ws.Range("D33").Value = "'501.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.89%  '
# Row 34
$ws.Range("E34").Value = '  +1.25%  '
# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '
# Row 36
$ws.Range("B36").Value = 'Cronos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D36").Value = "'0.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +22.22%  '
# Row 37
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = "'0.399"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +15.43%  '
# Row 38
$ws.Range("D38").Value = "'165.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.08%  '
# Row 39
$ws.Range("D39").Value = "'20.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.00%  '
# Row 40
$ws.Range("D40").Value = "'19.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.05%  '
# Row 41
$ws.Range("D41").Value = "'0.110"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.79%  '
# Row 42
$ws.Range("E42").Value = '  +0.06%  '
# Row 43
$ws.Range("D43").Value = "'180.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.60%  '
# Row 44
$ws.Range("D44").Value = "'4.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.47%  '
# Row 45
$ws.Range("D45").Value = "'1.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.87%  '
# Row 46
$ws.Range("E46").Value = '  -0.24%  '
# Row 47
$ws.Range("E47").Value = '  -3.66%  '
# Row 48
$ws.Range("D48").Value = "'0.592"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.72%  '
# Row 49
$ws.Range("D49").Value = "'3.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.04%  '
# Row 50
$ws.Range("D50").Value = "'2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.93%  '
# Row 51
$ws.Range("D51").Value = "'22.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.84%  '
